# Generate Report for Handoff
#
# The localization status report has progressed from "In Translation" to
# "Ready for handoff": the Status cells update, the handoff-generation
# timestamps advance a few seconds, and the Status/Datetime columns widen
# slightly (re-autofit) to accommodate the new, longer status text.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Status: "In Translation" -> "Ready for handoff" ------------------
# Overview!E2 (zh-cn status), Overview!F2 (de-de status)
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
# Status column on each language-specific sheet
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("C2").Value = "Ready for handoff"

# --- Latest HO Xliff Generate Date / Latest Handoff Datetime ----------
# Overview!G2 and de-de!H2 shared the same timestamp
$wsOverview.Range("G2").Value = "2016-08-27 19:05:32"
$wsDeDe.Range("H2").Value = "2016-08-27 19:05:32"
# zh-cn!H2 advanced to its own new handoff timestamp
$wsZhCn.Range("H2").Value = "2016-08-27 19:05:28"

# --- Column widths re-autofit for the now-wider Status/Datetime text --
# Target raw widths (~17.22 "characters") fall between the grid points
# this host's ColumnWidth setter can address; 16.33 lands on the closest
# reachable value.
$wsOverview.Columns.Item(5).ColumnWidth = 16.33
$wsOverview.Columns.Item(6).ColumnWidth = 16.33
$wsZhCn.Columns.Item(3).ColumnWidth = 16.33
$wsDeDe.Columns.Item(3).ColumnWidth = 16.33
